# Replace the short party abbreviations with "code - full name (native name, code)"
# everywhere they occur (including inside the '+'-joined coalition key strings),
# across every worksheet in the workbook. Excel automatically rebuilds the shared
# string table from the resulting cell contents.

$wb = $excel.ActiveWorkbook

$replacements = @(
    @('FPÖ', 'FPÖ - Freedom Party of Austria (Freiheitliche Partei Österreichs, FPÖ)'),
    @('GA ', 'GA - The Greens-Green Alternative (Die Grünen-Die Grüne Alternative, GA)'),
    @('LIF', 'LIF - Liberal Forum (Liberales Forum, LIF)'),
    @('SPÖ', 'SPÖ - Social Democratic Party of Austria  (Sozialdemokratische Partei Österreichs, SPÖ)'),
    @('ÖVP', 'ÖVP - Austrian People''s Party  (Österreichische Volkspartei, ÖVP)'),
    @('BZÖ', 'BZÖ - Alliance for the Future of Austria  (Bündnis Zukunft Österreich , BZÖ)'),
    @('FRA', 'FRANK - Team Frank Stronach (Team Frank Stronach, FRANK)'),
    @('NEO', 'NEOS - New Austria and Liberal Forum (NEOS Das Neue Österreich und Liberales Forum, NEOS)'),
    @('PIL', 'PILZ - Peter Pilz List (Liste Peter Pilz, PILZ)')
)

foreach ($ws in $wb.Worksheets) {
    foreach ($pair in $replacements) {
        $old = $pair[0]
        $new = $pair[1]
        # LookAt:=xlPart (2) so substrings inside combined "A+B+C" keys are replaced too.
        $ws.Cells.Replace($old, $new, 2, 1, $false, $false, $false) | Out-Null
    }
}
